$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells we touch to be stored as Text,
# so numeric-looking strings like "1.00" or "231.80" are not
# reinterpreted/rounded as numbers by Excel.
foreach ($addr in @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D8', 'D11', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D20', 'D22', 'D23', 'D25', 'D26', 'D27', 'D28', 'D30', 'D32', 'D33', 'D37', 'D38', 'D42', 'D43', 'D45', 'D48', 'D51')) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '27.338.30'
$ws.Cells.Item(2, 5).Value = '  -1.91%  '

$ws.Cells.Item(3, 4).Value = '1.655.63'
$ws.Cells.Item(3, 5).Value = '  -0.80%  '

$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.09%  '

$ws.Cells.Item(5, 4).Value = '213.19'
$ws.Cells.Item(5, 5).Value = '  -0.72%  '

$ws.Cells.Item(6, 4).Value = '0.514'
$ws.Cells.Item(6, 5).Value = '  -0.38%  '

$ws.Cells.Item(7, 4).Value = '1.00'
$ws.Cells.Item(7, 5).Value = '  -0.06%  '

$ws.Cells.Item(8, 4).Value = '23.65'
$ws.Cells.Item(8, 5).Value = '  +0.55%  '

$ws.Cells.Item(9, 5).Value = '  +0.05%  '

$ws.Cells.Item(10, 5).Value = '  -1.15%  '

$ws.Cells.Item(11, 4).Value = '0.0875'
$ws.Cells.Item(11, 5).Value = '  -0.42%  '

$ws.Cells.Item(12, 5).Value = '  -0.73%  '

$ws.Cells.Item(13, 4).Value = '1.662.75'
$ws.Cells.Item(13, 5).Value = '  -0.56%  '

$ws.Cells.Item(16, 4).Value = '65.78'
$ws.Cells.Item(16, 5).Value = '  -0.42%  '

$ws.Cells.Item(17, 4).Value = '27.353.74'
$ws.Cells.Item(17, 5).Value = '  -1.74%  '

$ws.Cells.Item(18, 4).Value = '231.80'
$ws.Cells.Item(18, 5).Value = '  -7.82%  '

$ws.Cells.Item(19, 5).Value = '  -0.87%  '

$ws.Cells.Item(20, 4).Value = '7.46'
$ws.Cells.Item(20, 5).Value = '  -1.09%  '

$ws.Cells.Item(21, 5).Value = '  -0.06%  '

$ws.Cells.Item(22, 4).Value = '4.37'
$ws.Cells.Item(22, 5).Value = '  -2.29%  '

$ws.Cells.Item(23, 4).Value = '9.36'
$ws.Cells.Item(23, 5).Value = '  +0.22%  '

$ws.Cells.Item(24, 5).Value = '  -1.72%  '

$ws.Cells.Item(25, 4).Value = '147.06'
$ws.Cells.Item(25, 5).Value = '  +0.12%  '

$ws.Cells.Item(26, 4).Value = '7.14'
$ws.Cells.Item(26, 5).Value = '  -1.18%  '

$ws.Cells.Item(27, 4).Value = '15.89'
$ws.Cells.Item(27, 5).Value = '  -2.75%  '

$ws.Cells.Item(28, 4).Value = '1.00'
$ws.Cells.Item(28, 5).Value = '  +0.05%  '

$ws.Cells.Item(29, 5).Value = '  -0.58%  '

$ws.Cells.Item(30, 4).Value = '0.0496'
$ws.Cells.Item(30, 5).Value = '  -1.06%  '

$ws.Cells.Item(31, 5).Value = '  -3.94%  '

$ws.Cells.Item(32, 4).Value = '3.30'
$ws.Cells.Item(32, 5).Value = '  -1.49%  '

$ws.Cells.Item(33, 4).Value = '1.449.76'
$ws.Cells.Item(33, 5).Value = '  +1.59%  '

$ws.Cells.Item(34, 5).Value = '  -0.56%  '

$ws.Cells.Item(35, 5).Value = '  +0.00%  '

$ws.Cells.Item(36, 5).Value = '  -0.74%  '

$ws.Cells.Item(37, 4).Value = '0.907'
$ws.Cells.Item(37, 5).Value = '  -2.49%  '

$ws.Cells.Item(38, 4).Value = '0.571'
$ws.Cells.Item(38, 5).Value = '  -2.18%  '

$ws.Cells.Item(39, 5).Value = '  -0.09%  '

$ws.Cells.Item(40, 5).Value = '  +0.47%  '

$ws.Cells.Item(41, 5).Value = '  -0.11%  '

$ws.Cells.Item(42, 4).Value = '5.50'
$ws.Cells.Item(42, 5).Value = '  +2.14%  '

$ws.Cells.Item(43, 4).Value = '65.08'
$ws.Cells.Item(43, 5).Value = '  -6.75%  '

$ws.Cells.Item(44, 5).Value = '  -0.27%  '

$ws.Cells.Item(45, 4).Value = '1.798.23'
$ws.Cells.Item(45, 5).Value = '  -0.69%  '

$ws.Cells.Item(46, 5).Value = '  -0.32%  '

$ws.Cells.Item(47, 5).Value = '  -0.81%  '

$ws.Cells.Item(48, 4).Value = '88.08'
$ws.Cells.Item(48, 5).Value = '  -1.11%  '

$ws.Cells.Item(49, 5).Value = '  -2.12%  '

$ws.Cells.Item(50, 5).Value = '  -0.46%  '

$ws.Cells.Item(51, 4).Value = '7.75'
$ws.Cells.Item(51, 5).Value = '  -0.95%  '

# Row 14/15: Polygon and Polkadot swap positions in the ranking
$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 4).Value = '4.08'
$ws.Cells.Item(14, 5).Value = '  -1.88%  '

$ws.Cells.Item(15, 2).Value = 'Polygon'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(15, 4).Value = '0.570'
$ws.Cells.Item(15, 5).Value = '  +3.39%  '

